$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 2323
    $ws.Range("F4").Value = 409
    $ws.Range("F6").Value = 6472
    $ws.Range("F7").Value = 324
}
